# Auto-generated Excel COM-interop script
# Applies crypto price/volume updates per commit:
# "Updated cryptos list on Sun Oct  8 20:34:32 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.906.47"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.635.49"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.34"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.257"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "1.867.02"
$ws.Range("D13").Value = "1.642.45"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "27.921.85"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("E19").Value = "  +2.74%  "
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").Value = "1.400.00"
$ws.Range("E35").Value = "  +3.05%  "
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.558"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.851"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("D46").Value = "1.775.75"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("E47").Value = "  -2.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.05%  "

Write-Output "Applied 70 cell updates"
